# Plantilla Lista de Tareas de la Entrega 6
# Implementación CU-30 y 31: se actualiza el estatus de varias tareas
# (de "Por inciar"/"Por iniciar" a "Hecho") y se registran horas
# consumidas en algunos dias para las tareas que ya fueron completadas.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")
$ws.Activate()

# --- Correccion ortografica del estatus inicial (fila 10 y 11) ---
$ws.Range("F10").Value = "Por iniciar"
$ws.Range("F11").Value = "Por iniciar"

# --- Tareas marcadas como terminadas ---
$ws.Range("F12").Value = "Hecho"
$ws.Range("F13").Value = "Hecho"
$ws.Range("F16").Value = "Hecho"
$ws.Range("F17").Value = "Hecho"
$ws.Range("F18").Value = "Hecho"
$ws.Range("F19").Value = "Hecho"

# --- Registro de horas consumidas por dia ---
$ws.Range("W12").Value = 0.5
$ws.Range("AF13").Value = 2
$ws.Range("Z16").Value = 1
$ws.Range("Z17").Value = 1
$ws.Range("W18").Value = 0.5
$ws.Range("AC19").Value = 1
$ws.Range("AF19").Value = 1

# Recalcular formulas dependientes (Rest. por dia y totales)
$excel.Calculate()

# Dejar la seleccion tal como quedo luego de registrar los datos
$ws.Range("X8").Select() | Out-Null
